# "Edit Projektplanung und Pflichtenheft"
# Update the actual-progress ("IST-" begin/duration) and "PROZENT ERLEDIGT"
# tracking columns (E:G) for each activity, and bump several "GEPLANTE(R)"
# begin dates (column C) to reflect the revised schedule.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Projekt")
if (-not $ws) { $ws = $wb.ActiveSheet }

# Row 9 - Planung
$ws.Range("E9").Value = 1
$ws.Range("F9").Value = 3
$ws.Range("G9").Value = 1

# Row 10 - Plichtenheft
$ws.Range("C10").Value = 4
$ws.Range("E10").Value = 4
$ws.Range("F10").Value = 3
$ws.Range("G10").Value = 1

# Row 11 - Testkatalog
$ws.Range("C11").Value = 7
$ws.Range("E11").Value = 7
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 1

# Row 12 - Programmierung
$ws.Range("C12").Value = 8
$ws.Range("E12").Value = 8
$ws.Range("F12").Value = 4
$ws.Range("G12").Value = 1

# Row 13 - Test
$ws.Range("C13").Value = 12
$ws.Range("E13").Value = 12
$ws.Range("F13").Value = 2
$ws.Range("G13").Value = 1

# Row 14 - Rückblende
$ws.Range("C14").Value = 14
$ws.Range("E14").Value = 14
$ws.Range("F14").Value = 2
$ws.Range("G14").Value = 1

# Nudge the "Period Highlight" spin-button control a fraction of a point
# (matches the sub-pixel reflow recorded when the sheet was last resaved).
try {
    $shp = $ws.Shapes.Item(1)
    $shp.Left = 936.7125196850394
    $shp.Top = 34.65
    $shp.Width = 53.0375
    $shp.Height = 18
} catch {
}

# Restore the author's cursor position at save time
$ws.Range("D9").Select()
